# Applies the "Notes" column (K) additions/updates for the predraft rookies sheet,
# matching the target commit ("Adds formatting for value"): new scouting-note text
# is written into column K for the listed rows, wrap text is enabled, and row heights
# are updated to fit the new note content.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$k16 = @'
Elic Ayomanor, a Canadian-born wide receiver with a unique high school path in the US, is entering the NFL draft after only two college seasons at Stanford. Despite a strong sophomore year with over 1,000 yards and a solid Breakout Score, his 2024 performance didn't significantly elevate his stats. He's a prototypical "X" receiver with a large frame and impressive 4.44 40-yard dash time, excelling in downfield targets. While his college numbers are good, not exceptional, his physical attributes and the demand for his position could lead to decent draft capital.
'@
$ws.Range("K16").Value = $k16
$ws.Range("K16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 101.5

$k18 = @'
The trend in football is shifting towards smaller, lighter wide receivers, but bigger players like Jayden Higgins stand out. At 6’4’’ and 214 pounds, Higgins combines size and smooth athleticism, which is rare for bigger receivers. He's versatile, playing from the slot and all over the formation.
Higgins' stats improved at Iowa State after transferring from Eastern Kentucky, showing solid production. His Breakout Score of 74.7 and strong yards per route run highlight his potential. Historically, players with similar profiles have had decent success in the NFL. Higgins’ blend of size, versatility, and performance makes him an intriguing prospect.
'@
$ws.Range("K18").Value = $k18
$ws.Range("K18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 145

$k21 = @'
Jaylin Noel, a versatile wide receiver from Iowa State, boasts above-average college stats, ranking in the 60-70th percentile among combine invitees. Primarily a slot receiver, he excelled in short-area targets and screens, evidenced by a low average depth of target (aDOT) early in his career. Despite this, he also demonstrated downfield capability, even surpassing his teammate Jayden Higgins in 20+ yard targets. Noel's versatility extends to special teams, where he returned kicks and punts. His strong combine performance has boosted his draft stock, with some models projecting him higher than others.
'@
$ws.Range("K21").Value = $k21
$ws.Range("K21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 116

$k24 = @'
Jalen Royals, after a brief stint at Georgia Military College, transferred to Utah State, where he had a slow start before a breakout junior season with 71 catches, 1,080 yards, and 15 touchdowns. His senior year further improved his receiving yards per team pass attempt. However, his overall Breakout Score of 43.9 is considered average, reflecting the level of competition at Utah State.
Projected as a Day 2 draft pick, Royals faces historical challenges, as Day 2 receivers with Breakout Scores below 50 have a lower success rate in producing high fantasy points. Despite this, he possesses athleticism (4.42 40-yard dash) and versatility, playing across the field. The ZAP Model compares him to Rashee Rice, though it acknowledges differences in their college usage and yards-after-catch statistics.
'@
$ws.Range("K24").Value = $k24
$ws.Range("K24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 159.5

$k26 = @'
Dylan Sampson, while possessing a 200-pound frame, doesn't align with the ideal size profile for high-end fantasy running backs. His production profile also doesn't compensate for this.
Across his three seasons at Tennessee, Sampson's rushing volume increased, but his receiving numbers remained underwhelming, culminating in a best-season reception share of only 8.5% and a below-average Breakout Score of 45.0.
His total yards per team play and yards after contact per attempt are also below average, and his tackle avoidance is considered ordinary. While he displayed some explosive runs, his absence at the combine leaves his speed untested.
Sampson's potential NFL workload is a concern. Historically, coaches tend to limit the usage of 200-pound backs, and his limited receiving ability further diminishes his fantasy prospects. Consequently, he risks a negative Draft Capital Delta.
'@
$ws.Range("K26").Value = $k26
$ws.Range("K26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 217.5

$k27 = @'
Jack Bech's NFL draft profile is intertwined with a recent personal tragedy: the loss of his brother, Tiger, in a New Orleans shooting. Bech honored his brother at the Senior Bowl, where he won MVP after scoring the game-winning touchdown. As a prospect, Bech had a promising freshman year at LSU, showing potential alongside future stars Brian Thomas Jr. and Malik Nabers. However, his role was primarily in the slot with short-range targets. After a quiet sophomore year, he transferred to TCU. His senior year was his most productive, but his overall analytical profile, including a below-average Breakout Score, is not particularly strong. Despite this, Bech's versatility, transitioning from a tight end/slot receiver to an "X" receiver, and his production against stiff competition at LSU, are noteworthy. While he has a favorable comparison to Michael Thomas, his overall draft projection is modest.
'@
$ws.Range("K27").Value = $k27
$ws.Range("K27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 159.5

$k28 = @'
Devin Neal's 4.58 40-yard dash time at the combine was underwhelming, failing to enhance his otherwise strong production profile. However, his overall ZAP Model score remains solid.
Neal demonstrated consistent production across his four seasons at Kansas, including consecutive 1,200-yard rushing campaigns. Notably, his receiving numbers were excellent, placing him among the top five running backs in this class with a best-season reception share of at least 13% and a Breakout Score above 85.
While his explosive run rate is below average, his pass protection skills, as graded by PFF, are a significant asset. This could facilitate early NFL playing time.
Despite his solid overall profile, the primary concern is a lack of explosive running ability, potentially limiting him to a handcuff role rather than a starting position. His success hinges on a team utilizing his pass protection and receiving abilities, similar to his top comparable, Shane Vereen.
'@
$ws.Range("K28").Value = $k28
$ws.Range("K28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 232

$k29 = @'
 Tory Horton looked pretty good through four years of college ball, and then he played in 
2024, his fifth campaign. He suffered a season-ending knee injury in October, totally 
derailing the season. Because he played that extra year and didn’t get any production out 
of it, the ZAP Model only saw that last season as an L. It did nothing but hurt his score.
 Back in 2020, Horton was at Nevada. He gave us a little production during his two years 
there, but he really started to cook when he transferred to Colorado State in 2022. That 
year -- his first season after transferring -- saw him reach a 3.51 prorated receiving 
yards per team pass attempt rate. That helped produce a Breakout Score of 70.4, the 11th
best mark in the class. Only three wide receiver combine invites had a season with as 
strong of a receiving yards per team pass attempt mark, too.
 Horton had a yards per route run against man that ranked 22nd in college football during 
his high-end 2022 season, and he can line up all over the field. He doesn’t have a thick 
build, but there’s good length. He can definitely play a role in the NFL.
 The question is, will he be good to go coming off of a season-ending injury? Well, this 
is where the NFL Combine can be helpful. Horton participated in the event, and he ran an 
impressive 4.41 40. He also did that at 196 pounds, when his listed weight in college was 
over 10 pounds below that.
 To me, that subjectively makes Horton even more impressive. He’s definitely a player to 
watch
'@
$ws.Range("K29").Value = $k29
$ws.Range("K29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 362.5

$k30 = @'
Xavier Restrepo, a slot receiver from Miami, lined up in the slot over 90% of his college routes, a high rate that's rare for top-100 draft picks. Historically, receivers with such high slot rates haven't found much NFL success, but the sample size is small.
Restrepo's production was limited until his final two seasons, where he achieved adjusted receiving yards per team pass attempt rates around 2.40, resulting in a Breakout Score of 48.9. While film analysts praise his route-running and spatial awareness, his production profile is less impressive compared to successful NFL slot receivers like Josh Downs, Christian Kirk, and Cooper Kupp, who had significantly higher Breakout Scores.
While Restrepo might find a role in the NFL, especially if he receives the projected draft capital, his ceiling as a fantasy football game-changer is questionable based on his college production.
'@
$ws.Range("K30").Value = $k30
$ws.Range("K30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 203

$k31 = @'
Ollie Gordon's 4.61 40-yard dash was mitigated by his significant 226-pound frame, resulting in a respectable Speed Score. He had a stellar 2023 season at Oklahoma State, leading the nation in rushing yards and winning the Doak Walker Award. 1  However, his production dipped in 2024.   
Gordon's size and downhill running style suggest he can handle a substantial workload. However, his lack of elite speed and below-average explosive play rate limit his upside, making volume crucial for fantasy success.
While players like James Conner and Rhamondre Stevenson have thrived with similar profiles, the depth of this running back class could lead to Gordon being part of a committee, limiting his immediate fantasy impact.
He has the potential for high volume in the right situation, but it's uncertain if he'll command a significant workload early in his NFL career.
'@
$ws.Range("K31").Value = $k31
$ws.Range("K31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 217.5

$k35 = @'
RJ Harvey presents an intriguing profile with a compact build similar to Jaylen Warren, but with a stronger production profile.
Strong Production: He boasts a solid 11.1% top-year reception share and an impressive 2.23 max-season total yards per team play rate, ranking highly in his class.
Combine Performance: His 4.40 40-yard dash resulted in a 109.4 Speed Score, highlighting his athleticism.1   
Age Concerns: At 24, Harvey's age is a significant factor. He began as a quarterback at Virginia in 2019, transferred to UCF, and faced setbacks including the pandemic and an ACL tear.
Historical Context: Historically, 24-year-old running back prospects have not yielded consistent, long-term fantasy success. Ray Davis and Tyrone Tracy are recent examples, but their NFL outcomes are yet to be determined.
Fantasy Outlook: Harvey's production and athleticism are appealing, but his age introduces risk. He could be a valuable contributor, but his long-term fantasy potential is uncertain
'@
$ws.Range("K35").Value = $k35
$ws.Range("K35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 217.5

$k36 = @'
 Jarquez Hunter has the third-best Breakout Score in this year’s draft class. It seems wild 
because he topped out at a 10.7% reception share at Auburn, which isn’t anything special. 
His Breakout Score is high, though, because, as a Sophomore, Hunter had 224 receiving yards 
on a team that threw it just 308 times. He was efficient with his touches.
 The combination of high-end Breakout Score and slightly-above-average reception share 
isn’t common. In fact, there’s only been one other running back in ZAP Model history with 
a Breakout Score above 90 (Hunter is at 93.8) and a best-season reception share below 11% 
(Hunter is at 10.7%). 
The other guy to do it? Kenneth Walker.
 Hunter didn’t have quite the same rushing performance as Walker in college, but don’t 
mistake that for him being a poor runner. Compared to the rest of the class, Hunter was 
fifth in career avoided tackles per rush, ninth in career explosive run rate, and ninth in 
career yards after contact per attempt. Shoutout to PFF for that data.
 Within a class littered with talent, Hunter stands out as a strong sleeper in the later 
rounds. He can do it all, and he showed at the NFL Combine that his speed is underrated, 
too, after running a 4.44 in the 40. Keep an eye on him throughout the draft process.
'@
$ws.Range("K36").Value = $k36
$ws.Range("K36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 348

$k38 = @'
Trevor Etienne is projected to have a "Low Risk" Draft Capital Delta due to his strong ZAP Model evaluation, primarily driven by his pass-catching abilities.
Despite his undersized 198-pound frame and limited collegiate workload (never exceeding 133 carries in a season), Etienne excels as a receiver. His top-season reception share ranks 10th in his class, and his Breakout Score is 11th.
However, his yards per team play rate never exceeded 1.25, reflecting his limited rushing volume. Notably, his best-season yards per team play is among the lowest in his class.
Despite this, historically, running backs with similar profiles—high Breakout Scores but low yards per team play—have found success as receivers in the NFL. Examples include Josh Jacobs, Tony Pollard, and James Cook.
Etienne's NFL role is likely as a pass-catching back, not a high-volume rusher. He's dynamic with the ball and projects as valuable PPR depth, even if he doesn't reach his brother's fantasy heights.
'@
$ws.Range("K38").Value = $k38
$ws.Range("K38").WrapText = $true
$ws.Rows.Item(38).RowHeight = 261

$k39 = @'
Ricky White boasts impressive raw production numbers, including high yards per route run and receiving yards per team pass attempt. However, these stats are inflated by the context of his competition at UNLV, in the Mountain West Conference.
After transferring from Michigan State, where he played behind NFL-caliber receivers, White excelled at UNLV. While his raw numbers are strong, the ZAP Model's adjustments for age and opponent strength result in a Breakout Score of 45.2. This creates a rare profile: a player with high raw production but a low adjusted score.
Historically, players with this profile have had mixed results, with John Brown being a notable success and Christian Watson still developing. This suggests that White's potential is uncertain.
As a fantasy analyst, the takeaway is to monitor White's draft capital and landing spot closely. While his raw numbers are enticing, the adjusted metrics suggest caution. He's a player with high potential but also significant risk, and his fantasy value will depend heavily on his NFL opportunity. Avoid overpaying for him based on his raw college stats alone.
'@
$ws.Range("K39").Value = $k39
$ws.Range("K39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 261

$k40 = @'
Savion Williams, a physically imposing 6'4", 220+ pound wide receiver from TCU, presents a unique profile due to his dual-threat ability as both a receiver and rusher. While his receiving production was notably low, with the worst adjusted receiving yards per team pass attempt rate in his draft class and a zero Breakout Score, his rushing stats are eye-catching. He recorded 51 carries for 322 yards and 6 touchdowns in his final college season. The ZAP Model values his versatility, particularly his rushing contribution, which boosts his overall evaluation. However, the model also acknowledges the significant risk associated with his poor receiving profile. Historically, wide receivers drafted in the top 100 without a Breakout Score have generally underperformed, with Terry McLaurin being a notable exception. Despite the analytical concerns, Williams's unique skill set could make him an intriguing Day 2 draft pick, offering a high ceiling if he lands in the right offensive system.
'@
$ws.Range("K40").Value = $k40
$ws.Range("K40").WrapText = $true
$ws.Rows.Item(40).RowHeight = 174

$k41 = @'
LeQuint Allen boasts an exceptional receiving profile, ranking near the top of his class in reception share and Breakout Score. However, his rushing metrics are concerning.
Receiving Dominance: He possesses elite receiving numbers, with a high reception share and Breakout Score, indicating significant pass-catching potential.
Rushing Concerns: His tackle avoidance and yards after contact per attempt are below average, raising questions about his effectiveness as a rusher.
Draft Capital Impact: His receiving profile is highly valuable, but it's heavily dependent on draft capital. A Day 3 selection would significantly diminish his fantasy prospects.
Historical Precedent: Historically, running backs with similar receiving profiles but lower draft capital have had limited success, with Kyren Williams being a notable exception.
Combine Uncertainty: Allen's decision to skip the combine leaves his speed and athleticism unknown, adding to the uncertainty surrounding his projection.
Fantasy Outlook: Despite the rushing concerns, Allen's pass-catching prowess makes him an intriguing PPR prospect, especially given his young age. His fantasy value is heavily tied to his draft position.
'@
$ws.Range("K41").Value = $k41
$ws.Range("K41").WrapText = $true
$ws.Rows.Item(41).RowHeight = 246.5

$k42 = @'
Tez Johnson, at 154 pounds, is one of the lightest wide receiver prospects in recent NFL draft history. While smaller receivers have found success in the NFL, Johnson's size, combined with a potentially poor combine performance, raises concerns about his draft capital and fantasy potential. Historically, only a few sub-170-pound receivers have achieved significant fantasy production, and those players generally had stronger production profiles and better draft capital than Johnson.
Johnson's college career spanned five years, with his most productive seasons at Oregon, where he posted impressive yards per route run numbers. However, his receiving yards per team pass attempt were more moderate.
While Johnson's yards per route run is a positive indicator, similar to recent successes like Ladd McConkey and Puka Nacua, his lack of experience playing outside the slot limits his versatility compared to other successful small receivers.
As a fantasy analyst, Johnson's size and limited perimeter experience suggest a lower ceiling. While he might find a role in the NFL, the odds of him becoming a high-end fantasy asset are lower than those of other small, but more versatile, receivers. Monitor his draft capital closely, but temper expectations.
'@
$ws.Range("K42").Value = $k42
$ws.Range("K42").WrapText = $true
$ws.Rows.Item(42).RowHeight = 232

$k43 = @'
Brashard Smith's unique journey from wide receiver at Miami to running back at SMU makes him an intriguing prospect.
* **Transition and Production:** He excelled in his single season as a running back, rushing for 1,332 yards and showcasing explosive speed with a sub-4.4 40-yard dash.
* **Comparison to Tyrone Tracy:** While similar in position change, Smith is smaller (194 pounds) than Tracy (209 pounds) and has less running back experience.
* **Size Concerns:** Historically, smaller running backs (under 200 pounds) rarely command high rushing workloads in the NFL.
* **Projected Role:** Smith is likely to begin his NFL career as a situational or mismatch player, leveraging his speed and agility.
* **Fantasy Outlook:** He's projected as a bench stash, with potential for more if he develops, but his size limits his likely early career volume.
'@
$ws.Range("K43").Value = $k43
$ws.Range("K43").WrapText = $true
$ws.Rows.Item(43).RowHeight = 188.5

$k44 = @'
We have a ton of great pass-catchers in this year’s draft class. So many of our 31 NFL 
Combine invites came with strong receiving numbers in the ZAP Model.
 Damien Martinez’s profile is lacking a bit there. His top year in prorated reception share 
is just 5.4%, second-worst in the class. He did have a 41.1 Breakout Score, but that’s not 
exactly what we’re looking for here.
 From 2011 to 2022, we’ve seen 44 running backs get drafted with best-season reception 
shares below 7%. Of those 44, only two scored more than 13 PPR points per game in one of 
their first three seasons.
 This isn’t everything for Martinez. After all, the backs who do succeed with his type 
of receiving numbers tend to be bigger in size, and Martinez was nearly 220 pounds at 
this year’s combine. He also ran a 4.51, giving him a 104.9 Speed Score, well above the 
threshold we’re looking for these guys to get to.
 We just probably won’t see Martinez deployed much as a receiver in the NFL. Going back to 
those 44 running backs with sub-7% best-season reception shares, only two of them were able 
to get to a double-digit percentage target share in one of their first three NFL seasons. 
And neither player was able to do it more than once. The vast majority of them were under 
the 5% target share mark as pros.
 Martinez is an angry runner who can lock up early-down work, but, especially in PPR 
formats, we’ll need more pass-catching for a true ceiling to be hit
'@
$ws.Range("K44").Value = $k44
$ws.Range("K44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 377

$k45 = @'
Tai Felton, despite his low BMI, excels at forcing missed tackles and generating yards after the catch, ranking highly in both categories. However, his overall production profile is average, with a top-season receiving yards per team pass attempt of 2.28 and a Breakout Score of 56.3.
Concerns exist about his play strength translating to the NFL, and his statistical comparisons are to players with similar builds who have struggled. His 4.37 speed is a potential asset, but relying solely on speed for NFL success is risky.
Overall, Felton's profile is unremarkable, leading the ZAP Model to project him as a benchwarmer. As a fantasy analyst, he presents a low-ceiling prospect with limited upside.
'@
$ws.Range("K45").Value = $k45
$ws.Range("K45").WrapText = $true
$ws.Rows.Item(45).RowHeight = 174

$k47 = @'
Montrell Johnson was faster than anticipated at the NFL Combine, running a 4.41 40-yard 
dash. At 212 pounds, that gave him the second-best Speed Score in the class, behind only 
Bhayshul Tuten. 
Should we care? I mean, kind of? It’s always nice when there’s a specific trait we can 
target with a late-round dart throw, and Speed Scores that high don’t grow on trees.
 Johnson’s performance at the combine was necessary for the ZAP Model, too, because there’s 
not a whole lot going for him outisde of it. His top season in total yards per team play 
was just 1.30 across his four years of college, fourth-worst in the class. When looking at 
the ZAP Model’s history, there’ve been just nine backs who scored 14 or more PPR points per 
game in either Year 1, Year 2, or Year 3 in the NFL who had a max-season total yards per 
team play rate below 1.40. Only four of them were drafted past Pick 100.
 One of those players is Isiah Pacheco, who’s Johnson’s top comp. Pacheco runs with more 
ferosity than almost any running back in the league, and that’s allowed him to find some 
success, but the two backs are similar in size, have elite Speed Scores, and both struggled 
a bit on the production front.
 Johnson could be a solid early-down runner with some upside on the right team in the NFL. 
We don’t get to say this a lot, but his combine performance saved him analytically
'@
$ws.Range("K47").Value = $k47
$ws.Range("K47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 391.5

$k48 = @'
 Pat Bryant’s stat lines across his four years at Illinois show pretty good yearly 
progression. He only had 98 yards as a Freshman, then it was 453 as a Sophomore, 560 as a 
Junior, and 984 as a Senior.
 But that Senior year -- 2024 -- looks much better in the ZAP Model than you might think. 
His prorated receiving yards per team pass attempt last year was 3.02, and, given the 
teams Illinois faced, his adjusted number looks even better. That’s more than double his 
receiving yards per team pass attempt mark from any other collegiate season.
 A Senior-year breakout isn’t something we’re desperately searching for, but we’re not 
talking about some high-end, obvious prospect who’s going to cost a lot during your rookie 
draft. Bryant came through with a solid 72.4 Breakout Score, something achieved by just 10 
wide receivers in this year’s class.
 From 2011 to 2022 –- this allows us to analyze wideouts in the ZAP Model database who’ve 
played three years in the NFL -- we’ve had 22 wide receivers hit a best-season yards per 
team pass attempt of 2.75 (Bryant was at 3.02) and a Breakout Score of 60 (Bryant was at 
72.4) while getting drafted after Pick 120. Of those 22 wideouts, 4 were able to score 
double-digit PPR points per game across their first three years in the league. That’s a 
rate of 18%.
 When looking at the group opposite of that, the 10-plus PPR points per game rate falls to 
4.9%. 
Late-round wide receivers don’t provide much fantasy relevancy. We get hits and remember 
when we do, but they’re rare. With Bryant, we’re at least looking at some marks that 
typically hit. Add on the fact that he’s a lengthy wide receiver who can play all over the 
formation, and you’ve got yourself a possible sleeper.
'@
$ws.Range("K48").Value = $k48
$ws.Range("K48").WrapText = $true
$ws.Rows.Item(48).RowHeight = 409.5

$k49 = @'
Nick Nash, a former quarterback turned wide receiver, had a remarkable 2024 season at San Jose State, winning the college football "triple crown" in receptions, receiving yards, and receiving touchdowns. 1  However, his raw production needs context.   
Despite his impressive raw numbers, his adjusted metrics, such as receiving yards per team pass attempt (2.71) and yards per route run (2.71), are merely average compared to other combine invitees. His age (almost 25) and the level of competition in the Mountain West Conference further temper his evaluation. His transition from quarterback to wide receiver is unique, but his Breakout Score of zero reflects the challenges posed by his age and program strength. While his 2024 season was undeniably productive, especially against larger programs, the raw numbers don't paint the full picture. As a fantasy analyst, Nash's size and recent success are intriguing, but his advanced age and adjusted metrics suggest caution. He's a developmental prospect with a potentially limited ceiling, and his raw stats should be viewed with skepticism.
'@
$ws.Range("K49").Value = $k49
$ws.Range("K49").WrapText = $true
$ws.Rows.Item(49).RowHeight = 203

# Update the sheet selection to match the saved cursor position after the edit.
[void]$ws.Range("K47").Select()

Write-Output "Done"